# Apply cryptos list price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.419.69"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.851.78"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'241.02"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'0.6295"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.07692"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("D10").Value = "'24.59"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.851.67"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "'5.030"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "'0.00001092"
$ws.Range("E14").Value = "  +8.76%  "
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.103.61"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'6.157"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "29.464.00"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'229.46"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'12.48"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'7.443"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("D25").Value = "'157.19"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'8.391"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'17.70"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "'1.316"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D31").Value = "'0.05719"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "'4.128"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "'4.055"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "'1.163"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'2.589"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'2.778"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'0.01793"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "1.221.18"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'6.494"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("D42").Value = "'0.9100"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "2.012.45"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'101.61"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'66.39"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "'7.141"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.012"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.683"
$ws.Range("E51").Value = "  +0.24%  "
